$d = $word.ActiveDocument

# The "Project Plan" table: row 5 (the "003" feature row) currently has its
# Bugfix columns (Bug #, Description, Finish date) empty. Fill in the Bug #
# and Description cells for a newly logged bug.
$table = $d.Tables(1)
$row = 5

$bugNumberCell = $table.Cell($row, 4)
$bugNumberCell.Range.Text = "003"

$bugDescCell = $table.Cell($row, 5)
$bugDescCell.Range.Text = "При регистрации ККТ при заполненных настройках ОФД в ККТ они не меняются на указанные в 1С"
